$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H88").Value = 1840083.5
$ws_ALC.Range("I88").Value = 2263.25
$ws_ALC.Range("J88").Value = 2452690.2
$ws_ALC.Range("K88").Value = 2263.25
$ws_ALC.Range("L88").Value = 2452690.2
$ws_ALC.Range("M88").Value = -1857.25
$ws_ALC.Range("N88").Value = -2453502.2
$ws_ALC.Range("H91").Value = 1840083.5
$ws_ALC.Range("I91").Value = 2263.25
$ws_ALC.Range("J91").Value = 2452690.2
$ws_ALC.Range("K91").Value = 2263.25
$ws_ALC.Range("L91").Value = 2452690.2
$ws_ALC.Range("M91").Value = -859.25
$ws_ALC.Range("N91").Value = -2455498.2
$ws_ALC.Range("H125").Value = 1390.6666
$ws_ALC.Range("I125").Value = 1011
$ws_ALC.Range("K125").Value = 9099
$ws_ALC.Range("M125").Value = -6639
$ws_ALC.Range("H135").Value = 2111.42
$ws_ALC.Range("I135").Value = 983.8372000000001
$ws_ALC.Range("J135").Value = 9038
$ws_ALC.Range("K135").Value = 8854.534800000001
$ws_ALC.Range("L135").Value = 81342
$ws_ALC.Range("M135").Value = -6319.534800000001
$ws_ALC.Range("N135").Value = -86412
$ws_ALC.Range("H138").Value = 3591.254
$ws_ALC.Range("I138").Value = 2075.1936
$ws_ALC.Range("J138").Value = 5059.9375
$ws_ALC.Range("K138").Value = 6225.5808
$ws_ALC.Range("L138").Value = 15179.8125
$ws_ALC.Range("M138").Value = -1085.5808
$ws_ALC.Range("N138").Value = -25459.8125
$ws_ALC.Range("H141").Value = 2317
$ws_ALC.Range("I141").Value = 2225.484
$ws_ALC.Range("J141").Value = 3026.25
$ws_ALC.Range("K141").Value = 6676.451999999999
$ws_ALC.Range("L141").Value = 9078.75
$ws_ALC.Range("M141").Value = -1496.451999999999
$ws_ALC.Range("N141").Value = -19438.75

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H132").Value = 1096.1522
$ws_ARM.Range("I132").Value = 866.2432
$ws_ARM.Range("K132").Value = 2598.7296
$ws_ARM.Range("M132").Value = -68.72960000000012

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H118").Value = 44712
$ws_BSM.Range("J118").Value = 44712
$ws_BSM.Range("L118").Value = 44712
$ws_BSM.Range("N118").Value = -48026
$ws_BSM.Range("H123").Value = 50000
$ws_BSM.Range("J123").Value = 50000
$ws_BSM.Range("L123").Value = 50000
$ws_BSM.Range("N123").Value = -59800
$ws_BSM.Range("H127").Value = 0
$ws_BSM.Range("J127").Value = 0
$ws_BSM.Range("L127").Value = 0
$ws_BSM.Range("N127").ClearContents()
$ws_BSM.Range("H134").Value = 3844.2952
$ws_BSM.Range("I134").Value = 1268.2778
$ws_BSM.Range("K134").Value = 3804.8334
$ws_BSM.Range("M134").Value = -1269.8334

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 5437111
$ws_CRP.Range("I31").Value = 7247767.5
$ws_CRP.Range("J31").Value = 5141
$ws_CRP.Range("K31").Value = 7247767.5
$ws_CRP.Range("L31").Value = 5141
$ws_CRP.Range("M31").Value = -7247472.5
$ws_CRP.Range("N31").Value = -5731
$ws_CRP.Range("H34").Value = 5437111
$ws_CRP.Range("I34").Value = 7247767.5
$ws_CRP.Range("J34").Value = 5141
$ws_CRP.Range("K34").Value = 7247767.5
$ws_CRP.Range("L34").Value = 5141
$ws_CRP.Range("M34").Value = -7247565.5
$ws_CRP.Range("N34").Value = -5545
$ws_CRP.Range("H58").Value = 1130.174
$ws_CRP.Range("I58").Value = 945.23334
$ws_CRP.Range("J58").Value = 1476.9375
$ws_CRP.Range("K58").Value = 945.23334
$ws_CRP.Range("L58").Value = 1476.9375
$ws_CRP.Range("M58").Value = -742.23334
$ws_CRP.Range("N58").Value = -1882.9375
$ws_CRP.Range("H86").Value = 37513.5
$ws_CRP.Range("I86").Value = 45475
$ws_CRP.Range("J86").Value = 29552
$ws_CRP.Range("K86").Value = 45475
$ws_CRP.Range("L86").Value = 29552
$ws_CRP.Range("M86").Value = -44352
$ws_CRP.Range("N86").Value = -31798
$ws_CRP.Range("H89").Value = 37513.5
$ws_CRP.Range("I89").Value = 45475
$ws_CRP.Range("J89").Value = 29552
$ws_CRP.Range("K89").Value = 227375
$ws_CRP.Range("L89").Value = 147760
$ws_CRP.Range("M89").Value = -221759
$ws_CRP.Range("N89").Value = -158992
$ws_CRP.Range("H132").Value = 1187.362
$ws_CRP.Range("I132").Value = 1077.5098
$ws_CRP.Range("J132").Value = 1987.7142
$ws_CRP.Range("K132").Value = 3232.5294
$ws_CRP.Range("L132").Value = 5963.142599999999
$ws_CRP.Range("M132").Value = -702.5294000000004
$ws_CRP.Range("N132").Value = -11023.1426
$ws_CRP.Range("H136").Value = 1130.174
$ws_CRP.Range("I136").Value = 945.23334
$ws_CRP.Range("J136").Value = 1476.9375
$ws_CRP.Range("K136").Value = 2835.70002
$ws_CRP.Range("L136").Value = 4430.8125
$ws_CRP.Range("M136").Value = -285.7000200000002
$ws_CRP.Range("N136").Value = -9530.8125

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 870.02563
$ws_CUL.Range("I5").Value = 567.0909
$ws_CUL.Range("K5").Value = 1701.2727
$ws_CUL.Range("M5").Value = -1589.2727
$ws_CUL.Range("H12").Value = 15.045455
$ws_CUL.Range("I12").Value = 1
$ws_CUL.Range("J12").Value = 17.263159
$ws_CUL.Range("K12").Value = 3
$ws_CUL.Range("L12").Value = 51.78947700000001
$ws_CUL.Range("M12").Value = 170
$ws_CUL.Range("N12").Value = -397.789477
$ws_CUL.Range("H38").Value = 1103.6
$ws_CUL.Range("I38").Value = 879.3333
$ws_CUL.Range("J38").Value = 2000.6666
$ws_CUL.Range("K38").Value = 2637.9999
$ws_CUL.Range("L38").Value = 6001.9998
$ws_CUL.Range("M38").Value = -2290.9999
$ws_CUL.Range("N38").Value = -6695.9998
$ws_CUL.Range("H92").Value = 498.57144
$ws_CUL.Range("I92").Value = 498
$ws_CUL.Range("K92").Value = 1494
$ws_CUL.Range("M92").Value = -246
$ws_CUL.Range("H131").Value = 10000840
$ws_CUL.Range("I131").Value = 766.4286
$ws_CUL.Range("J131").Value = 13889758
$ws_CUL.Range("K131").Value = 2299.2858
$ws_CUL.Range("L131").Value = 41669274
$ws_CUL.Range("M131").Value = 2740.7142
$ws_CUL.Range("N131").Value = -41679354
$ws_CUL.Range("H132").Value = 1251.5641
$ws_CUL.Range("I132").Value = 1221.3158
$ws_CUL.Range("J132").Value = 1280.3
$ws_CUL.Range("K132").Value = 10991.8422
$ws_CUL.Range("L132").Value = 11522.7
$ws_CUL.Range("M132").Value = -8461.842200000001
$ws_CUL.Range("N132").Value = -16582.7
$ws_CUL.Range("H135").Value = 870.02563
$ws_CUL.Range("I135").Value = 567.0909
$ws_CUL.Range("K135").Value = 5103.8181
$ws_CUL.Range("M135").Value = -2568.8181

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H132").Value = 1631.0176
$ws_GSM.Range("I132").Value = 1621.9387
$ws_GSM.Range("J132").Value = 1686.625
$ws_GSM.Range("K132").Value = 4865.8161
$ws_GSM.Range("L132").Value = 5059.875
$ws_GSM.Range("M132").Value = -2335.8161
$ws_GSM.Range("N132").Value = -10119.875

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H132").Value = 2371.9443
$ws_LTW.Range("I132").Value = 2567.6487
$ws_LTW.Range("J132").Value = 1946
$ws_LTW.Range("K132").Value = 7702.946100000001
$ws_LTW.Range("L132").Value = 5838
$ws_LTW.Range("M132").Value = -5172.946100000001
$ws_LTW.Range("N132").Value = -10898

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H136").Value = 856.68085
$ws_WVR.Range("I136").Value = 837.1429000000001
$ws_WVR.Range("J136").Value = 913.6667
$ws_WVR.Range("K136").Value = 2511.4287
$ws_WVR.Range("L136").Value = 2741.0001
$ws_WVR.Range("M136").Value = 38.57129999999961
$ws_WVR.Range("N136").Value = -7841.0001
